$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.855756402015686
$ws.Range("B1").Value = 3.144791126251221
$ws.Range("C1").Value = 3.184075593948364
$ws.Range("D1").Value = 3.625637769699097
$ws.Range("E1").Value = 3.137698173522949
